# feat: add 2022-Q1 data
#
# - Inserts a new worksheet "2022-Q1" (fund-holding detail, same shape as
#   "2021-Q4") positioned between "2021-Q4" and "总计".
# - Adds a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data (existing "2021-Q4" row shifts down).

$wb = $excel.ActiveWorkbook
$wsQ4   = $wb.Worksheets.Item("2021-Q4")
$wsTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before "总计" and copy the
#    header / index-column formatting from the "2021-Q4" sheet so the
#    new sheet matches the existing look (bold + bordered + centered).
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($wsTotal)
$wsQ1.Name = "2022-Q1"

# Re-resolve "总计" by name: inserting a sheet can leave earlier
# worksheet variables bound to the wrong tab position.
$wsTotal = $wb.Worksheets.Item("总计")

$wsQ4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$wsQ4.Range("A2").Copy()
$wsQ1.Range("A2:A7").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------
$wsQ1.Range("B1").Value = "基金代码"
$wsQ1.Range("C1").Value = "基金名称"
$wsQ1.Range("D1").Value = "基金规模"
$wsQ1.Range("E1").Value = "股票总仓位"
$wsQ1.Range("F1").Value = "仓位占比"
$wsQ1.Range("G1").Value = "持有市值(亿元)"
$wsQ1.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 3. Data rows. Columns B-G hold numeric-looking text (fund codes /
#    percentages) that must stay text (leading zeros, fixed decimals),
#    so those columns are pre-formatted as Text before being written.
#    Column H (rank) is a genuine number, column A is the numeric index.
# ---------------------------------------------------------------------
$wsQ1.Range("B2:B7").NumberFormat = "@"
$wsQ1.Range("D2:G7").NumberFormat = "@"

$data = @(
    @("159869", "华夏中证动漫游戏ETF",       "6.20", "98.75", "3.37", "0.2089", 10),
    @("516010", "国泰中证动漫游戏ETF",       "4.95", "98.91", "3.29", "0.1629", 10),
    @("001628", "招商体育文化休闲股票",       "2.95", "83.21", "4.59", "0.1354", 4),
    @("516770", "华泰柏瑞中证动漫游戏ETF",   "1.11", "96.56", "3.32", "0.0369", 10),
    @("011231", "光大保德信锦弘混合A",       "4.13", "20.96", "0.76", "0.0314", 1),
    @("011232", "光大保德信锦弘混合C",       "1.29", "20.96", "0.76", "0.0098", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $wsQ1.Cells.Item($r, 1).Value = $i
    $wsQ1.Cells.Item($r, 2).Value = $row[0]
    $wsQ1.Cells.Item($r, 3).Value = $row[1]
    $wsQ1.Cells.Item($r, 4).Value = $row[2]
    $wsQ1.Cells.Item($r, 5).Value = $row[3]
    $wsQ1.Cells.Item($r, 6).Value = $row[4]
    $wsQ1.Cells.Item($r, 7).Value = $row[5]
    $wsQ1.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 4. "总计" sheet: insert a new data row above the existing "2021-Q4"
#    row and populate it with the 2022-Q1 totals. Clear the formatting
#    Excel auto-inherits on row insert so the new plain-data cells stay
#    unstyled (matching the rest of the sheet's data rows).
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("A2:D2").ClearFormats()

$wsQ4.Range("A2").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)   # re-apply bold/border index style

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 6
$wsTotal.Range("D2").Value = 0.59

$wsTotal.Range("A3").Value = 1
